$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.776.26'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '3.729.52'
$ws.Range("E3").Value = '  -2.32%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.14%  '

$ws.Range("D7").Value = '3.728.46'
$ws.Range("E7").Value = '  -2.36%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("E10").Value = '  +2.35%  '

$ws.Range("E11").Value = '  +1.83%  '

$ws.Range("E12").Value = '  -1.41%  '

$ws.Range("E13").Value = '  -1.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000246'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("D15").Value = '4.354.35'
$ws.Range("E15").Value = '  -1.85%  '

$ws.Range("D16").Value = '3.725.36'
$ws.Range("E16").Value = '  -2.10%  '

$ws.Range("D17").Value = '68.778.42'
$ws.Range("E17").Value = '  +1.46%  '

$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '493.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +15.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.726'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000144'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.87%  '

$ws.Range("E26").Value = '  -3.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("E28").Value = '  -1.36%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.61%  '

$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("E32").Value = '  +1.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.30%  '

$ws.Range("D34").Value = '3.871.88'
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("E35").Value = '  -1.15%  '

$ws.Range("D36").Value = '3.661.60'
$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.21%  '

$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '439.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.36%  '

$ws.Range("E43").Value = '  -0.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.23%  '

$ws.Range("E45").Value = '  +0.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.41%  '

$ws.Range("D50").Value = '2.776.15'
$ws.Range("E50").Value = '  -2.58%  '

$ws.Range("E51").Value = '  +0.55%  '

